$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add Wins / Losses / Ties headers in AD1:AF1, matching the
# bold/centered/bordered style used by the other header cells (copy from AC1).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats (values already set above)

# Data rows 2-49: team record columns (same value repeated for every player row)
for ($r = 2; $r -le 49; $r++) {
    $ws.Cells.Item($r, 30).Value = 88   # AD
    $ws.Cells.Item($r, 31).Value = 74   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
